# adding gas flow equation. not yet working
#
# This script reproduces, cell-by-cell, the target state of the "edge" and
# "device" sheets described by the commit's diff:
#  - edge: nodeFrom/nodeTo rename, capacity bumps for the electrical edges,
#          a new el edge (node1 -> node1ex), dropping the old "flow_k"
#          column in favour of a trailing pressureFrom/pressureTo/gasflow_k
#          block for the gas edges, and a brand-new gas edge row.
#  - device: header rename ("comment" -> "name") + in_*/out_* column swap,
#            several device renames, and assorted capacity/pressure value
#            edits (incl. the 1e10 "infinite" Pmax values for the new
#            compressor/separator devices).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "edge"
# ---------------------------------------------------------------------
$edge = $wb.Worksheets.Item("edge")

# Header row (A1:L1) -- rewrite in full since a column (flow_k) was
# dropped and three columns (pressureFrom/pressureTo/gasflow_k) were
# appended after "comment".
$edge.Range("A1:L1").ClearContents()
$edge.Cells.Item(1,1).Value  = "type"
$edge.Cells.Item(1,2).Value  = "nodeFrom"
$edge.Cells.Item(1,3).Value  = "nodeTo"
$edge.Cells.Item(1,4).Value  = "include"
$edge.Cells.Item(1,5).Value  = "capacity"
$edge.Cells.Item(1,6).Value  = "reistance"
$edge.Cells.Item(1,7).Value  = "reactance"
$edge.Cells.Item(1,8).Value  = "distance"
$edge.Cells.Item(1,9).Value  = "comment"
$edge.Cells.Item(1,10).Value = "pressureFrom"
$edge.Cells.Item(1,11).Value = "pressureTo"
$edge.Cells.Item(1,12).Value = "gasflow_k"

# Data rows 2:9 -- clear the old 2:8 block first (row 5 moves content
# down to row 6, a new row 5 is inserted, and a new row 9 is appended),
# then write the full new table.
$edge.Range("A2:L9").ClearContents()

# el edges (capacity raised 100 -> 200)
$edge.Cells.Item(2,1).Value = "el"
$edge.Cells.Item(2,2).Value = "node1"
$edge.Cells.Item(2,3).Value = "node2"
$edge.Cells.Item(2,4).Value = 1
$edge.Cells.Item(2,5).Value = 200
$edge.Cells.Item(2,6).Value = 0.001
$edge.Cells.Item(2,7).Value = 0.01
$edge.Cells.Item(2,8).Value = 2

$edge.Cells.Item(3,1).Value = "el"
$edge.Cells.Item(3,2).Value = "node1"
$edge.Cells.Item(3,3).Value = "node3"
$edge.Cells.Item(3,4).Value = 1
$edge.Cells.Item(3,5).Value = 200
$edge.Cells.Item(3,6).Value = 0.001
$edge.Cells.Item(3,7).Value = 0.01
$edge.Cells.Item(3,8).Value = 2

$edge.Cells.Item(4,1).Value = "el"
$edge.Cells.Item(4,2).Value = "node2"
$edge.Cells.Item(4,3).Value = "node3"
$edge.Cells.Item(4,4).Value = 1
$edge.Cells.Item(4,5).Value = 200
$edge.Cells.Item(4,6).Value = 0.001
$edge.Cells.Item(4,7).Value = 0.01
$edge.Cells.Item(4,8).Value = 3

# new el export edge
$edge.Cells.Item(5,1).Value = "el"
$edge.Cells.Item(5,2).Value = "node1"
$edge.Cells.Item(5,3).Value = "node1ex"
$edge.Cells.Item(5,4).Value = 1
$edge.Cells.Item(5,5).Value = 200

# gas edges -- now carry pressureFrom/pressureTo/gasflow_k (J/K/L)
$edge.Cells.Item(6,1).Value = "gas"
$edge.Cells.Item(6,2).Value = "node4"
$edge.Cells.Item(6,3).Value = "node3"
$edge.Cells.Item(6,4).Value = 1
$edge.Cells.Item(6,5).Value = 100
$edge.Cells.Item(6,8).Value = 1
$edge.Cells.Item(6,10).Value = 100
$edge.Cells.Item(6,11).Value = 95
$edge.Cells.Item(6,12).Value = 0.001

$edge.Cells.Item(7,1).Value = "gas"
$edge.Cells.Item(7,2).Value = "node3"
$edge.Cells.Item(7,3).Value = "node2"
$edge.Cells.Item(7,4).Value = 1
$edge.Cells.Item(7,5).Value = 100
$edge.Cells.Item(7,8).Value = 1
$edge.Cells.Item(7,10).Value = 100
$edge.Cells.Item(7,11).Value = 400
$edge.Cells.Item(7,12).Value = 0.001

$edge.Cells.Item(8,1).Value = "gas"
$edge.Cells.Item(8,2).Value = "node2"
$edge.Cells.Item(8,3).Value = "node1"
$edge.Cells.Item(8,4).Value = 1
$edge.Cells.Item(8,5).Value = 100
$edge.Cells.Item(8,8).Value = 1
$edge.Cells.Item(8,10).Value = 400
$edge.Cells.Item(8,11).Value = 395
$edge.Cells.Item(8,12).Value = 0.001

$edge.Cells.Item(9,1).Value = "gas"
$edge.Cells.Item(9,2).Value = "node1"
$edge.Cells.Item(9,3).Value = "node1ex"
$edge.Cells.Item(9,4).Value = 1
$edge.Cells.Item(9,5).Value = 200
$edge.Cells.Item(9,8).Value = 1
$edge.Cells.Item(9,10).Value = 400
$edge.Cells.Item(9,11).Value = 395
$edge.Cells.Item(9,12).Value = 0.001

$edge.Activate()
$edge.Range("J10").Select()

# ---------------------------------------------------------------------
# Sheet "device"
# ---------------------------------------------------------------------
$device = $wb.Worksheets.Item("device")

# Header rename + in_*/out_* column swap
$device.Cells.Item(1,2).Value = "name"
$device.Cells.Item(1,5).Value  = "in_el"
$device.Cells.Item(1,6).Value  = "in_gas"
$device.Cells.Item(1,7).Value  = "in_heat"
$device.Cells.Item(1,8).Value  = "out_el"
$device.Cells.Item(1,9).Value  = "out_gas"
$device.Cells.Item(1,10).Value = "out_heat"

# Row 3 (gas turbine): in_el 0 instead of -0.45, + out_el 0.45
$device.Cells.Item(3,5).Value = 0
$device.Cells.Item(3,8).Value = 0.45

# Row 4 ("gas processing" -> "compressor"), Pmax -> 1e10 (sci format)
$device.Cells.Item(4,2).Value = "compressor"
$device.Cells.Item(4,14).Value = 10000000000
$device.Cells.Item(4,14).NumberFormat = "0.00E+00"

# Row 5 ("gas compressor (gas)" -> "compressor1"), in_el 0.1 -> 0, Pmax -> 1e10
$device.Cells.Item(5,2).Value = "compressor1"
$device.Cells.Item(5,5).Value = 0
$device.Cells.Item(5,14).Value = 10000000000
$device.Cells.Item(5,14).NumberFormat = "0.00E+00"

# Row 6 ("gas compressor" -> "separator"), Pmax -> 1e10
$device.Cells.Item(6,2).Value = "separator"
$device.Cells.Item(6,14).Value = 10000000000
$device.Cells.Item(6,14).NumberFormat = "0.00E+00"

# Row 8 ("wind generator" -> "wind turb"), Pmax 100 -> 400, Pmin 0 -> 300
$device.Cells.Item(8,2).Value = "wind turb"
$device.Cells.Item(8,14).Value = 400
$device.Cells.Item(8,15).Value = 300

# Row 12 (wellhead): Pmin 200 -> 250
$device.Cells.Item(12,15).Value = 250

# Row 13 ("gas export" -> "gas EXPORT"), Pmax 200 -> 300
$device.Cells.Item(13,2).Value = "gas EXPORT"
$device.Cells.Item(13,14).Value = 300

# Row 14 ("el export" -> "el EXPORT"), node moved to node1ex
$device.Cells.Item(14,1).Value = "node1ex"
$device.Cells.Item(14,2).Value = "el EXPORT"

# Row 15 ("diesel aggregator backup" -> "diesel backup")
$device.Cells.Item(15,2).Value = "diesel backup"

$device.Activate()
$device.Range("N12").Select()
